# Rever_DailyTrack_BALRAJ_2022.xlsx - "Add files via upload"
#
# Populates rows 14-15 of the MAR-22 sheet with the 5th daily-task entry
# (previously blank placeholder rows):
#   Row 14: No=5, Date=07-Mar-2022, App=RPA GSS,   Comment about Service
#           Order Pending stale-error issue, Status=WIP
#   Row 15: (same task) App=RPA RLOGIC, Comment about the MFA fix at BLR,
#           % complete=90%, Status=WIP

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MAR-22")

$xlPasteFormats = -4122

# ---- 1. Copy the cell formatting from the nearest fully-populated rows so
#         the new cells inherit the same borders / number formats / wrap
#         alignment used throughout the sheet, rather than staying with the
#         bare "empty placeholder" style. -------------------------------

# Row 14 formats (mirrors row 2's layout: No / Date / App / wrapped Task / Status)
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A14").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("B2").Copy() | Out-Null
$ws.Range("B14").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("C2").Copy() | Out-Null
$ws.Range("C14").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("D2").Copy() | Out-Null
$ws.Range("D14").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("F2").Copy() | Out-Null
$ws.Range("F14").PasteSpecial($xlPasteFormats) | Out-Null

# Row 15 formats (App / wrapped Task / % complete / Status, continuing the task)
$ws.Range("C2").Copy() | Out-Null
$ws.Range("C15").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("D4").Copy() | Out-Null
$ws.Range("D15").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("E2").Copy() | Out-Null
$ws.Range("E15").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("F2").Copy() | Out-Null
$ws.Range("F15").PasteSpecial($xlPasteFormats) | Out-Null

$excel.CutCopyMode = $false

# ---- 2. Fill in the values --------------------------------------------

$ws.Range("A14").Value = 5
$ws.Range("B14").Value = 44627
$ws.Range("C14").Value = "RPA GSS"
$ws.Range("D14").Value = "1.   The task of Service Order Pending  is work in progress, and clicking at each cell is success whereas click at each row is having `nstale error issue. It is a challenging  for us to download."
$ws.Range("F14").Value = "WIP"

$ws.Range("C15").Value = "RPA RLOGIC"
$ws.Range("D15").Value = "1. MFA issue has been fixed at BLR center only, and Following that the daily tasks are success, and we have downloaded the Feb22 files for the GL and PL files are going on for the BLR."
$ws.Range("E15").Value = 0.9
$ws.Range("F15").Value = "WIP"

# ---- 3. Row heights so the wrapped comment text is fully visible -------

$ws.Rows.Item(14).RowHeight = 43.2
$ws.Rows.Item(15).RowHeight = 28.8

# ---- 4. Move the active selection to D15, matching where the author's
#         cursor ended up after typing the new comment. -----------------

$ws.Range("D15").Select() | Out-Null
